$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.817127704620361
$ws.Range("B1").Value = 2.614512920379639
$ws.Range("C1").Value = 2.834985971450806
$ws.Range("D1").Value = 3.4032883644104
$ws.Range("E1").Value = 0.9723404049873352
